$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextCell "D2" '51.781.76'
$ws.Range("E2").Value = '  -0.71%  '

Set-TextCell "D3" '2.916.53'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("E4").Value = '  -0.07%  '

Set-TextCell "D5" '355.89'
$ws.Range("E5").Value = '  +1.36%  '

Set-TextCell "D6" '109.20'
$ws.Range("E6").Value = '  -2.54%  '

Set-TextCell "D7" '0.565'
$ws.Range("E7").Value = '  +1.46%  '

$ws.Range("E8").Value = '  +0.04%  '

Set-TextCell "D9" '0.624'
$ws.Range("E9").Value = '  -1.08%  '

Set-TextCell "D10" '39.04'
$ws.Range("E10").Value = '  -2.08%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell "D11" '0.137'
$ws.Range("E11").Value = '  +1.17%  '

$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell "D12" '0.0868'
$ws.Range("E12").Value = '  +0.79%  '

Set-TextCell "D13" '19.51'
$ws.Range("E13").Value = '  -2.15%  '

Set-TextCell "D14" '7.79'
$ws.Range("E14").Value = '  -0.11%  '

Set-TextCell "D15" '3.377.78'
$ws.Range("E15").Value = '  +0.40%  '

Set-TextCell "D16" '2.925.35'
$ws.Range("E16").Value = '  +0.40%  '

Set-TextCell "D17" '0.980'
$ws.Range("E17").Value = '  -1.28%  '

Set-TextCell "D18" '51.813.28'
$ws.Range("E18").Value = '  -0.77%  '

$ws.Range("E19").Value = '  +2.11%  '

$ws.Range("E20").Value = '  -1.20%  '

Set-TextCell "D21" '13.91'
$ws.Range("E21").Value = '  -1.81%  '

Set-TextCell "D22" '0.0₃0977'
$ws.Range("E22").Value = '  -0.09%  '

Set-TextCell "D23" '70.49'
$ws.Range("E23").Value = '  -0.38%  '

Set-TextCell "D24" '268.65'
$ws.Range("E24").Value = '  -0.34%  '

$ws.Range("E25").Value = '  +0.91%  '

$ws.Range("E26").Value = '  +13.49%  '

Set-TextCell "D27" '7.69'
$ws.Range("E27").Value = '  +18.97%  '

Set-TextCell "D28" '26.88'
$ws.Range("E28").Value = '  +0.70%  '

$ws.Range("E29").Value = '  +0.16%  '

$ws.Range("E30").Value = '  +11.13%  '

Set-TextCell "D31" '10.50'
$ws.Range("E31").Value = '  -0.69%  '

Set-TextCell "D32" '37.41'
$ws.Range("E32").Value = '  -0.29%  '

Set-TextCell "D33" '6.05'
$ws.Range("E33").Value = '  -0.71%  '

Set-TextCell "D34" '52.19'
$ws.Range("E34").Value = '  -1.64%  '

$ws.Range("E35").Value = '  -8.72%  '

Set-TextCell "D36" '0.0441'
$ws.Range("E36").Value = '  -2.11%  '

Set-TextCell "D37" '1.00'
$ws.Range("E37").Value = '  -0.05%  '

$ws.Range("E38").Value = '  +3.84%  '

Set-TextCell "D39" '18.22'
$ws.Range("E39").Value = '  -2.24%  '

$ws.Range("E40").Value = '  -3.39%  '

Set-TextCell "D41" '2.72'
$ws.Range("E41").Value = '  -3.86%  '

$ws.Range("E42").Value = '  +2.32%  '

Set-TextCell "D43" '22.75'
$ws.Range("E43").Value = '  -2.68%  '

Set-TextCell "D44" '118.99'
$ws.Range("E44").Value = '  -1.61%  '

$ws.Range("E45").Value = '  -0.09%  '

$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell "D46" '2.48'
$ws.Range("E46").Value = '  -5.01%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell "D47" '3.46'
$ws.Range("E47").Value = '  -1.67%  '

Set-TextCell "D48" '2.120.97'
$ws.Range("E48").Value = '  -3.51%  '

Set-TextCell "D49" '0.248'
$ws.Range("E49").Value = '  -4.72%  '

Set-TextCell "D50" '0.0339'
$ws.Range("E50").Value = '  +1.37%  '

$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell "D51" '9.10'
$ws.Range("E51").Value = '  +0.46%  '
